$d = $word.ActiveDocument

# 1. Merge the "2023 " and "(lots of overlap here!)" runs into a single run.
$d.Content.Find.Execute("2023 (lots of overlap here!)", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2023 (lots of overlap here!)", 2)

# 2. Add a new bullet after "2023 (~30%)" noting low tow coverage years.
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.ListFormat.ListLevelNumber = 1
$newPara.Range.InsertAfter("Tow coverage is low in 2018, 2021 and 2022")
